# Apply the "check record" feature-doc edit to Sheet1 (row 16) plus the
# small cosmetic changes (row height, selection, column widths) that came
# along with it in the author's commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16 content: replace the old "add remark" text with the new
#     "check record" text -------------------------------------------------
$ws.Range("B16").Value = "检测记录"
$ws.Range("C16").Value = "添加检测记录"
$ws.Range("D16").Value = "点击添加备注图标，进入检测记录screen，支持用户录入文档，点击确定后提交文档，并返回检测详情screen"

# --- Row 16 grew very slightly taller to fit the new wrapped text --------
$ws.Rows.Item(16).RowHeight = 17.15

# --- Columns got very slightly narrower (re-flow after the text edit) ----
$ws.Columns.Item(1).ColumnWidth = 7.417
$ws.Columns.Item(2).ColumnWidth = 7.417
$ws.Columns.Item(3).ColumnWidth = 41.584
$ws.Columns.Item(4).ColumnWidth = 91.25

# --- Active cell moved from D14 to D16 while the author was editing it ---
$ws.Range("D16").Select()

# --- Tab scroll ratio nudged slightly in the book view --------------------
$win = $excel.ActiveWindow
$win.TabRatio = 0.995
